# Remove NP10 from short names
# The NP is in the directory and it's not meaningful to most ppl
#
# Column F ("short_name") on sheet "all_runs" has "(NP10)" / "(NP10, ...)"
# qualifiers stripped out (the NP## is already captured in the "directory"
# column, so repeating it in the short name is redundant). The previously
# plain "Pathway 3a" / "Pathway 3b" / "Pathway 4 - No New Pricing" short
# names (from the NP07 runs) are disambiguated with "(NP07)" now that the
# bare names are freed up by the NP10 rows losing their suffix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_runs")

# Pathway 1a
$ws.Range("F32").Value = "Pathway 1a"
$ws.Range("F33").Value = "Pathway 1a (more VZ, less trn, lower toll)"

# Pathway 1b
$ws.Range("F52").Value = "Pathway 1b"
$ws.Range("F53").Value = "Pathway 1b (more VZ, less trn, lower toll)"

# Pathway 2a
$ws.Range("F60").Value = "Pathway 2a with 20pct art toll"
$ws.Range("F61").Value = "Pathway 2a with 10pct art toll"
$ws.Range("F62").Value = "Pathway 2a with 20pct art toll (more VZ, less trn, lower toll)"
$ws.Range("F63").Value = "Pathway 2a with 10pct art toll (more VZ, less trn, lower toll)"

# Pathway 2b
$ws.Range("F71").Value = "Pathway 2b with 20pct art toll"
$ws.Range("F72").Value = "Pathway 2b with 10pct art toll"
$ws.Range("F73").Value = "Pathway 2b with 20pct art toll (more VZ, less trn, lower toll)"
$ws.Range("F74").Value = "Pathway 2b with 10pct art toll (more VZ, less trn, lower toll)"

# Pathway 3a (older NP07 rows now disambiguated; NP10 rows lose the suffix)
$ws.Range("F79").Value = "Pathway 3a (NP07)"
$ws.Range("F80").Value = "Pathway 3a (NP07)"
$ws.Range("F84").Value = "Pathway 3a"
$ws.Range("F85").Value = "Pathway 3a (more VZ, less trn)"

# Pathway 3b
$ws.Range("F86").Value = "Pathway 3b (NP07)"
$ws.Range("F87").Value = "Pathway 3b (NP07)"
$ws.Range("F91").Value = "Pathway 3b"
$ws.Range("F92").Value = "Pathway 3b (more VZ, less trn)"

# Pathway 4
$ws.Range("F93").Value = "Pathway 4 - No New Pricing (NP07)"
$ws.Range("F97").Value = "Pathway 4 - No New Pricing"
$ws.Range("F98").Value = "Pathway 4 - No New Pricing (more VZ, less trn)"

# View state: freeze just the header row (ySplit=1) but scroll the frozen
# pane down near the bottom of the data, and leave the selection on F96.
$ws.Activate() | Out-Null
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F96").Select() | Out-Null
